# Update "paises.xlsx" - countries & provincias Spain
# 1) Swap the ranking order of "Republica del Chad" and "Suazilandia"
#    (Suazilandia's updated figures moved it ahead of Republica del Chad)
# 2) Refresh the COVID-19 statistics for a number of countries
# 3) Update the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names between row 151 and row 152 ---------------------
$name151 = $ws.Range("A151").Value2
$name152 = $ws.Range("A152").Value2
$ws.Range("A151").Value = $name152
$ws.Range("A152").Value = $name151

# --- Refresh numeric data (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) -----------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1216820
$ws.Range("C4").Value = 3985
$ws.Range("D4").Value = 188205
$ws.Range("E4").Value = 958028
$ws.Range("G4").Value = 666
$ws.Range("H4").Value = 70587

# Row 6 - Italia
$ws.Range("B6").Value = 213013
$ws.Range("C6").Value = 1075
$ws.Range("D6").Value = 85231
$ws.Range("E6").Value = 98467
$ws.Range("F6").Value = 1427
$ws.Range("G6").Value = 236
$ws.Range("H6").Value = 29315

# Row 9 - Alemania
$ws.Range("B9").Value = 166424
$ws.Range("C9").Value = 272
$ws.Range("E9").Value = 24331

# Row 15 - Canada
$ws.Range("B15").Value = 61165
$ws.Range("C15").Value = 393
$ws.Range("D15").Value = 26305
$ws.Range("E15").Value = 30942
$ws.Range("G15").Value = 64
$ws.Range("H15").Value = 3918

# Row 29 - Singapur
$ws.Range("D29").Value = 1519
$ws.Range("E29").Value = 17873
$ws.Range("F29").Value = 24

# Row 36 - Polonia
$ws.Range("B36").Value = 14431
$ws.Range("C36").Value = 425
$ws.Range("E36").Value = 9435
$ws.Range("G36").Value = 18
$ws.Range("H36").Value = 716

# Row 47 - Noruega
$ws.Range("B47").Value = 7928
$ws.Range("C47").Value = 24
$ws.Range("E47").Value = 7681

# Row 69 - Grecia
$ws.Range("B69").Value = 2642
$ws.Range("C69").Value = 10
$ws.Range("E69").Value = 1122

# Row 100 - Sri Lanka
$ws.Range("B100").Value = 762
$ws.Range("C100").Value = 11
$ws.Range("D100").Value = 213
$ws.Range("E100").Value = 540

# Row 120 - Jordania
$ws.Range("B120").Value = 471
$ws.Range("C120").Value = 6
$ws.Range("D120").Value = 375
$ws.Range("E120").Value = 87

# Row 151 - now "Suazilandia" (updated data)
$ws.Range("B151").Value = 119
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 12
$ws.Range("E151").Value = 106
$ws.Range("H151").Value = 1

# Row 152 - now "Republica del Chad" (retains previous data)
$ws.Range("B152").Value = 117
$ws.Range("D152").Value = 39
$ws.Range("E152").Value = 68
$ws.Range("H152").Value = 10

# --- Update the "last updated" timestamp text ----------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 18:08"
